$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-10
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13)
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
